# Applies the "fixed extract users from teams and removed park visualizer" edit:
#  - Replaces the sample/demo rows on the "Create Teams" sheet (old Spanish/French
#    example row) with the new Italian BKW rows, adding three extra rows.
#  - Switches which sheet/selection is active: "Create Teams" becomes the
#    selected tab (was "NAS Downloads"), with A2:E5 selected; "NAS Downloads"
#    keeps its own remembered selection (B14) but is no longer the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Create Teams")
$ws2 = $wb.Worksheets.Item("NAS Downloads")

# --- Replace / extend the data rows on "Create Teams" -----------------------
$rows = @(
    @("0-IT-CTL-01", "EBBK0001", "ZP1", "I801", "BKW"),
    @("0-IT-CTL-02", "EBBK0001", "ZP1", "I801", "BKW"),
    @("0-IT-VLC-01",  "EBBK0001", "ZP1", "I801", "804I"),
    @("0-IT-BZI-01",  "ETBK0001", "ZP1", "I810", "BKW")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $ws1.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# --- Update remembered selections / active tab -------------------------------
# Give "NAS Downloads" its new remembered selection first (it loses focus).
$ws2.Range("B14").Select()

# Make "Create Teams" the active sheet/tab, with A2:E5 selected.
$ws1.Select()
$ws1.Range("A2:E5").Select()
